$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.78000000000043
$ws.Range("G2").Value = [double]"6.079547976156618e-10"
$ws.Range("H2").Value = [double]"2.830620666353716e-09"
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = 39.84822472011997
$ws.Range("L2").Value = "[27.513550208699826, 52.182899231540105]"
$ws.Range("M2").Value = [double]"1.384267589443766e-09"
$ws.Range("N2").Value = [double]"2.768535178887532e-09"
$ws.Range("O2").Value = 1.515763422452732
$ws.Range("P2").Value = "[1.1258159859711165, 1.9057108589343477]"
$ws.Range("Q2").Value = [double]"8.968381592922015e-13"
$ws.Range("R2").Value = [double]"8.968381592922015e-13"
$ws.Range("S2").Value = 63.85532348432839
$ws.Range("T2").Value = "[55.731968111292275, 71.97867885736451]"
$ws.Range("W2").Value = 18.80204204204237
$ws.Range("X2").Value = 17.26414414414445
$ws.Range("Y2").Value = 20.3399399399403

# Row 3 updates
$ws.Range("E3").Value = 25.40000000000053
$ws.Range("G3").Value = [double]"1.060539434050156e-10"
$ws.Range("H3").Value = [double]"1.65668924362116e-09"
$ws.Range("K3").Value = 37.59990811139389
$ws.Range("L3").Value = "[25.468456274557923, 49.731359948229866]"
$ws.Range("M3").Value = [double]"4.1496874825242e-09"
$ws.Range("N3").Value = [double]"4.1496874825242e-09"
$ws.Range("O3").Value = 1.792500312859042
$ws.Range("P3").Value = "[1.4277107755052736, 2.157289850212811]"
$ws.Range("S3").Value = 63.78156880818275
$ws.Range("T3").Value = "[56.488898270417835, 71.07423934594766]"
$ws.Range("W3").Value = 18.15375375375413
$ws.Range("X3").Value = 16.67907907907943
$ws.Range("Y3").Value = 19.62842842842883
